$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 344-345; existing rows 344:365 shift down to 346:367
$ws.Rows("344:345").Insert()

# Fill in the new row 344
$ws.Range("A344").Value = 6
$ws.Range("B344").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C344").Value = 'Metropolitana'
$ws.Range("D344").Value = 44578
$ws.Range("E344").Value = 13
$ws.Range("F344").Value = 100112030
$ws.Range("G344").Value = 'Poroto granado'
$ws.Range("H344").Value = 'Sin especificar'
$ws.Range("I344").Value = 'Primera'
$ws.Range("J344").Value = 340
$ws.Range("K344").Value = 25000
$ws.Range("L344").Value = 28000
$ws.Range("M344").Value = 26412
$ws.Range("N344").Value = '$/saco 25 kilos'
$ws.Range("O344").Value = 'Región Metropolitana'
$ws.Range("P344").Value = 1056
$ws.Range("Q344").Value = 25
$ws.Range("R344").Value = 'Hortaliza'

# Fill in the new row 345
$ws.Range("A345").Value = 6
$ws.Range("B345").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C345").Value = 'Metropolitana'
$ws.Range("D345").Value = 44578
$ws.Range("E345").Value = 13
$ws.Range("F345").Value = 100112030
$ws.Range("G345").Value = 'Poroto granado'
$ws.Range("H345").Value = 'Sin especificar'
$ws.Range("I345").Value = 'Primera'
$ws.Range("J345").Value = 370
$ws.Range("K345").Value = 27000
$ws.Range("L345").Value = 30000
$ws.Range("M345").Value = 28216
$ws.Range("N345").Value = '$/saco 25 kilos'
$ws.Range("O345").Value = 'Región del Maule'
$ws.Range("P345").Value = 1129
$ws.Range("Q345").Value = 25
$ws.Range("R345").Value = 'Hortaliza'
